# Adds the numeric-to-categorical "Orange" results block to the
# "Numeric To Categorical" sheet, mirroring the existing Python block,
# and updates sheet selection / active tab state to match.

$wb = $excel.ActiveWorkbook

$wsNum = $wb.Worksheets.Item("Numeric To Categorical")
$wsOut = $wb.Worksheets.Item("Outlier Removal")
$wsMiss = $wb.Worksheets.Item("Missing Values")

# --- Restructure "Numeric To Categorical": shift existing Python block
# down one row and right one column, to make room for a header row and
# an "Orange" results block (matching the layout already used on the
# "Outlier Removal" sheet). ---
$wsNum.Rows("1:1").Insert()
$wsNum.Columns("A:A").Insert()

# Column D should use the wrap-text variant of the column C font/style.
$wsNum.Range("D2:D6").WrapText = $true

# --- New header row (Python / Orange / Data Polish) ---
$wsNum.Range("C1").Value2 = "Python"
$wsNum.Range("G1").Value2 = "Orange"
$wsNum.Range("K1").Value2 = "Data Polish"

# Copy header-row formatting (s=14 / s=15) from the equivalent cells on
# the "Outlier Removal" sheet, which already has this exact header.
$wsOut.Range("C1").Copy()
$wsNum.Range("C1").PasteSpecial(-4122)
$wsOut.Range("D1").Copy()
$wsNum.Range("D1").PasteSpecial(-4122)
$wsOut.Range("G1").Copy()
$wsNum.Range("G1").PasteSpecial(-4122)
$wsOut.Range("H1").Copy()
$wsNum.Range("H1").PasteSpecial(-4122)
$wsOut.Range("K1").Copy()
$wsNum.Range("K1").PasteSpecial(-4122)

# --- New "Orange" block (columns F:H), mirroring the Python block ---
$wsNum.Range("F2").Value2 = "Action"
$wsNum.Range("G2").Value2 = "Time"
$wsNum.Range("H2").Value2 = "Content"

$wsNum.Range("F3").Value2 = "Load Data"
$wsNum.Range("G3").Value2 = "2 min"
$wsNum.Range("H3").Value2 = "Use ""File"" widget to load the dataset"

$wsNum.Range("F4").Value2 = "Convert Types"
$wsNum.Range("G4").Value2 = "3 min"
$wsNum.Range("H4").Value2 = "Use ""Edit Domain"" widget to change attribute type"

$wsNum.Range("F5").Value2 = "Verify Changes"
$wsNum.Range("G5").Value2 = "1 min"
$wsNum.Range("H5").Value2 = "Use ""Data Table"" widget to verify changes"

$wsNum.Range("F6").Value2 = "Overall"
$wsNum.Range("G6").Value2 = "6 min"

# Copy formatting for the Orange block + row-2 header (F:H) from the
# equivalent, already-styled cells on "Outlier Removal".
$wsOut.Range("F2:H2").Copy()
$wsNum.Range("F2:H2").PasteSpecial(-4122)

$wsOut.Range("F3:H3").Copy()
$wsNum.Range("F3:H3").PasteSpecial(-4122)
$wsNum.Range("F4:H4").PasteSpecial(-4122)
$wsNum.Range("F5:H5").PasteSpecial(-4122)

$wsOut.Range("F7:H7").Copy()
$wsNum.Range("F6:H6").PasteSpecial(-4122)

# Column widths to match the new content (mirrors "Outlier Removal").
$wsNum.Columns("B:B").ColumnWidth = 19
$wsNum.Columns("F:F").ColumnWidth = 17

# Row heights for the newly-sized rows.
$wsNum.Rows("1:1").RowHeight = 19
$wsNum.Rows("2:2").RowHeight = 18
$wsNum.Rows("3:3").RowHeight = 90
$wsNum.Rows("4:4").RowHeight = 108
$wsNum.Rows("5:5").RowHeight = 108
$wsNum.Rows("6:6").RowHeight = 54
$wsNum.Rows("7:7").RowHeight = 17

# --- Sheet-view / selection bookkeeping ---
# "Missing Values" was the active tab; "Numeric To Categorical" becomes
# the active tab instead.
$wsOut.Range("A1:XFD1").Select()
$wsOut.Range("C1").Activate()

$wsNum.Activate()
$wsNum.Range("J4").Select()
